# Sprint 3 Hours Log -- Mason
# Log a new entry (row 7) for SF-13: the tile-removal button.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 42845
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "SF-13"
$ws.Range("D7").Value = "Added a button that will remove all tiles of a number based on user input"

# Leave the selection where the user would end up after finishing the row.
$ws.Range("A8").Select()
